$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("G10").Value = 1.83
$ws.Range("I10").Value = 4.75
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("AI10").Value = 21
$ws.Range("AQ10").Value = 41

# Row 11
$ws.Range("G11").Value = 1.42
$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 2
$ws.Range("L11").Value = 8
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 3.25
$ws.Range("Q11").Value = 2.08
$ws.Range("R11").Value = 1.73
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 6
$ws.Range("Z11").Value = 9
$ws.Range("AA11").Value = 13
$ws.Range("AB11").Value = 34
$ws.Range("AD11").Value = 8.5
$ws.Range("AE11").Value = 26
$ws.Range("AJ11").Value = 23
$ws.Range("AS11").Value = 201
$ws.Range("AU11").Value = 10
$ws.Range("AZ11").Value = 201
$ws.Range("BA11").Value = 251
$ws.Range("BD11").Value = 151

# Row 15
$ws.Range("G15").Value = 2.4
$ws.Range("I15").Value = 2.9
$ws.Range("J15").Value = 3.1
$ws.Range("L15").Value = 3.5
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.3
$ws.Range("P15").Value = 3.4
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.75
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 1.91
$ws.Range("W15").Value = 8
$ws.Range("X15").Value = 11
$ws.Range("Y15").Value = 9.5
$ws.Range("Z15").Value = 23
$ws.Range("AA15").Value = 21
$ws.Range("AB15").Value = 29
$ws.Range("AC15").Value = 9.5
$ws.Range("AG15").Value = 251
$ws.Range("AH15").Value = 9
$ws.Range("AI15").Value = 15
$ws.Range("AJ15").Value = 11
$ws.Range("AK15").Value = 29
$ws.Range("AL15").Value = 23
$ws.Range("AM15").Value = 34
$ws.Range("AN15").Value = 4.5
$ws.Range("AO15").Value = 13
$ws.Range("AP15").Value = 23
$ws.Range("AQ15").Value = 41
$ws.Range("AR15").Value = 67
$ws.Range("AS15").Value = 151
$ws.Range("AW15").Value = 4.75
$ws.Range("AX15").Value = 17
$ws.Range("AY15").Value = 26
$ws.Range("AZ15").Value = 51
$ws.Range("BA15").Value = 81
$ws.Range("BB15").Value = 201

# Row 17
$ws.Range("G17").Value = 3.1
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 2.2
$ws.Range("J17").Value = 3.6
$ws.Range("L17").Value = 2.88
$ws.Range("N17").Value = 13
$ws.Range("Q17").Value = 1.75
$ws.Range("R17").Value = 2.05
$ws.Range("AH17").Value = 9
$ws.Range("AN17").Value = 5

# Row 18
$ws.Range("G18").Value = 1.67
$ws.Range("I18").Value = 5.25
$ws.Range("J18").Value = 2.3
$ws.Range("L18").Value = 5.5
$ws.Range("M18").Value = 1.06
$ws.Range("N18").Value = 10
$ws.Range("U18").Value = 1.91
$ws.Range("V18").Value = 1.91
$ws.Range("X18").Value = 7.5
$ws.Range("Z18").Value = 13
$ws.Range("AA18").Value = 13
$ws.Range("AH18").Value = 13
$ws.Range("AI18").Value = 26
$ws.Range("AJ18").Value = 17
$ws.Range("AK18").Value = 51
$ws.Range("AL18").Value = 41
$ws.Range("AN18").Value = 3.6
$ws.Range("AO18").Value = 8.5
$ws.Range("AQ18").Value = 29
$ws.Range("AU18").Value = 8.5
$ws.Range("AW18").Value = 6.5
$ws.Range("AX18").Value = 26
$ws.Range("AY18").Value = 34
$ws.Range("AZ18").Value = 101
$ws.Range("BA18").Value = 126
$ws.Range("BB18").Value = 251

# Row 19
$ws.Range("G19").Value = 2.45
$ws.Range("H19").Value = 3.1
$ws.Range("I19").Value = 3.1
$ws.Range("J19").Value = 3.2
$ws.Range("K19").Value = 2.05
$ws.Range("L19").Value = 3.6
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 8.5
$ws.Range("O19").Value = 1.36
$ws.Range("P19").Value = 3.2
$ws.Range("S19").Value = 1.44
$ws.Range("T19").Value = 2.63
$ws.Range("U19").Value = 1.91
$ws.Range("V19").Value = 1.91
$ws.Range("W19").Value = 7.5
$ws.Range("Y19").Value = 10
$ws.Range("Z19").Value = 23
$ws.Range("AC19").Value = 8.5
$ws.Range("AH19").Value = 8.5
$ws.Range("AJ19").Value = 11
$ws.Range("AL19").Value = 26
$ws.Range("AM19").Value = 34
$ws.Range("AO19").Value = 15
$ws.Range("AQ19").Value = 51
$ws.Range("AT19").Value = 2.63
$ws.Range("AU19").Value = 8
$ws.Range("AV19").Value = 51
$ws.Range("AW19").Value = 4.75
$ws.Range("AX19").Value = 17
$ws.Range("AZ19").Value = 51
$ws.Range("BB19").Value = 201

# Row 20
$ws.Range("G20").Value = 2.6
$ws.Range("I20").Value = 2.9
$ws.Range("J20").Value = 3.2
$ws.Range("K20").Value = 2.05
$ws.Range("L20").Value = 3.5
$ws.Range("Q20").Value = 2.08
$ws.Range("R20").Value = 1.73
$ws.Range("W20").Value = 8
$ws.Range("X20").Value = 12
$ws.Range("Y20").Value = 10
$ws.Range("Z20").Value = 23
$ws.Range("AA20").Value = 21
$ws.Range("AB20").Value = 29
$ws.Range("AE20").Value = 15
$ws.Range("AH20").Value = 8.5
$ws.Range("AI20").Value = 13
$ws.Range("AJ20").Value = 11
$ws.Range("AK20").Value = 29
$ws.Range("AL20").Value = 23
$ws.Range("AM20").Value = 34
$ws.Range("AN20").Value = 4.5
$ws.Range("AO20").Value = 15
$ws.Range("AP20").Value = 23
$ws.Range("AR20").Value = 67
$ws.Range("AW20").Value = 4.75
$ws.Range("AX20").Value = 17
$ws.Range("AY20").Value = 26
$ws.Range("BA20").Value = 81
$ws.Range("BB20").Value = 201

# Row 21
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 1.95
$ws.Range("L21").Value = 2.63
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 9.5
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.82
$ws.Range("U21").Value = 1.8
$ws.Range("V21").Value = 1.95
$ws.Range("W21").Value = 11
$ws.Range("AA21").Value = 34
$ws.Range("AC21").Value = 9.5
$ws.Range("AD21").Value = 6.5
$ws.Range("AG21").Value = 251
$ws.Range("AI21").Value = 9
$ws.Range("AK21").Value = 17
$ws.Range("AS21").Value = 201
$ws.Range("AX21").Value = 11
$ws.Range("AY21").Value = 21
$ws.Range("BA21").Value = 51

# Row 22
$ws.Range("G22").Value = 2.7
$ws.Range("H22").Value = 3.25
$ws.Range("I22").Value = 2.63
$ws.Range("L22").Value = 3.25
$ws.Range("Q22").Value = 2.02
$ws.Range("R22").Value = 1.88
$ws.Range("W22").Value = 9
$ws.Range("Z22").Value = 26

# Row 26
$ws.Range("G26").Value = 2.63
$ws.Range("I26").Value = 3
$ws.Range("M26").Value = 1.13
$ws.Range("N26").Value = 6
$ws.Range("W26").Value = 6.5
$ws.Range("X26").Value = 11
$ws.Range("AC26").Value = 6
$ws.Range("AP26").Value = 34
$ws.Range("AV26").Value = 81
$ws.Range("AX26").Value = 19
$ws.Range("AZ26").Value = 67
$ws.Range("BB26").Value = 351

# Row 27
$ws.Range("G27").Value = 2.3
$ws.Range("I27").Value = 3.25
$ws.Range("J27").Value = 3.1
$ws.Range("AB27").Value = 34
$ws.Range("AO27").Value = 13
